$d = $word.ActiveDocument

# 1. Title line: "디자인 팀의 주요 업무 문서" -> "팀 핵심 책임 문서 디자인"
#    Also the run's bold toggles on (w:b w:val="0" -> w:b, i.e. bold=true).
$r = $d.Content.Find.Execute("디자인 팀의 주요 업무 문서", $true, $false, $false, $false, $false,
                              $true, 1, $false, "팀 핵심 책임 문서 디자인", 2)

# Find the run we just replaced and turn bold on.
$rng = $d.Content
$rng.Find.Execute("팀 핵심 책임 문서 디자인", $true, $false, $false, $false, $false,
                   $true, 1, $false, "", 0)
$rng.Font.Bold = 1

# 2. Purpose sentence
$d.Content.Find.Execute(": 이 문서에서는 그래픽 디자인 연구소의 모든 디자인 팀 구성원의 핵심 책임을 간략하게 설명합니다.", $true, $false, $false, $false, $false,
                          $true, 1, $false, ": 이 문서에서는 Graphic Design Institute의 모든 디자인 팀 구성원의 핵심 책임을 간략하게 설명합니다.", 2)

# 3. Collaboration responsibility text
$d.Content.Find.Execute(": 다른 디자이너, 개발자 및 이해 관계자와 협력하여 프로젝트 요구 사항을 충족하는 고품질 디자인을 만듭니다. ", $true, $false, $false, $false, $false,
                          $true, 1, $false, ": 다른 디자이너, 개발자, 이해 관계자들과 공동으로 작업을 진행하여 프로젝트 요구 사항을 충족하는 고품질 디자인을 만들어야 합니다. ", 2)

# 4. Design responsibility text
$d.Content.Find.Execute(": 사용자에게 친숙하고 접근성이 뛰어나며 반응성이 뛰어난 시각적으로 매력적인 디자인을 만듭니다. ", $true, $false, $false, $false, $false,
                          $true, 1, $false, ": 사용자들이 익숙한 방식으로 쉽게 활용할 수 있으며 적극적으로 반응할 수 있는 멋진 스타일의 디자인을 만들어야 합니다. ", 2)

# 5. Communication responsibility text
$d.Content.Find.Execute(": 팀 구성원, 이해 관계자 및 클라이언트와 효과적으로 통신하여 프로젝트 요구 사항을 충족하는지 확인합니다. ", $true, $false, $false, $false, $false,
                          $true, 1, $false, ": 프로젝트 요구 사항을 충족할 수 있도록 팀 구성원, 이해 관계자, 고객과 효율적으로 커뮤니케이션합니다. ", 2)

# 6. Research responsibility text (the word "연구" inside this sentence also changes to "리서치")
$d.Content.Find.Execute(": 디자인 결정을 알리기 위해 사용자 요구 사항, 기본 설정 및 동작을 식별하기 위한 연구를 수행합니다. ", $true, $false, $false, $false, $false,
                          $true, 1, $false, ": 디자인 결정을 알리기 위해 사용자 요구 사항, 기본 설정 및 동작을 식별하기 위한 리서치를 수행합니다. ", 2)

# 7. "연구" heading run -> "리서치" (the standalone run that starts the same paragraph as #6).
#    Find with MatchWholeWord is unreliable on CJK text (no word-boundary concept), and a
#    plain Find would also hit the "연구" inside "연구소" above, so locate the run directly:
#    it is the very first two characters of the paragraph that now reads "연구: 디자인 결정을...".
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t.StartsWith("연구:")) {
        $start = $p.Range.Start
        $rng = $d.Range($start, $start + 2)
        $rng.Text = "리서치"
        break
    }
}

# 8. Usability test responsibility text
$d.Content.Find.Execute(": 디자인이 사용자 요구를 충족하고 모든 사용자가 액세스할 수 있도록 유용성 테스트를 수행합니다. ", $true, $false, $false, $false, $false,
                          $true, 1, $false, ": 사용 편의성 테스트를 수행하여 디자인이 사용자의 요구를 충족하며 모든 사용자가 쉽게 사용 가능한 상태인지를 확인해야 합니다. ", 2)

# 9. "전문 개발" -> "전문적인 개발"
$d.Content.Find.Execute("전문 개발", $true, $false, $false, $false, $false,
                          $true, 1, $false, "전문적인 개발", 2)
